$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete row 552 entirely (the post "「あす日が昇るだろう。私は朝が大好きだ」").
# This shifts all subsequent rows (553..717) up by one, matching the diff
# (which shows rows 553-717 renumbered to 552-716 with identical content,
# and the sheet dimension shrinking from A1:C717 to A1:C716).
$ws.Rows.Item(552).Delete()
